$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B2 = New-Object "object[,]" 24,5
$arr_B2[0,0] = 2.157254401178534
$arr_B2[0,1] = 0.4977145929672133
$arr_B2[0,2] = 0.03186391388105481
$arr_B2[0,3] = 0.07750110998557247
$arr_B2[0,4] = 3.371396017997341
$arr_B2[1,0] = 2.039924398498954
$arr_B2[1,1] = 0.4649510196448716
$arr_B2[1,2] = 0.03163220842753223
$arr_B2[1,3] = 0.07768356930806775
$arr_B2[1,4] = 3.329408434567185
$arr_B2[2,0] = 1.96927304524803
$arr_B2[2,1] = 0.4451547636298585
$arr_B2[2,2] = 0.03148807221397121
$arr_B2[2,3] = 0.07781798545151464
$arr_B2[2,4] = 3.3055279485515
$arr_B2[3,0] = 1.940828943214342
$arr_B2[3,1] = 0.4371674501138614
$arr_B2[3,2] = 0.03142884699622872
$arr_B2[3,3] = 0.0778784010526552
$arr_B2[3,4] = 3.296271856290517
$arr_B2[4,0] = 1.936126699630393
$arr_B2[4,1] = 0.4358459600357776
$arr_B2[4,2] = 0.03141898264304288
$arr_B2[4,3] = 0.07788877398071126
$arr_B2[4,4] = 3.294763533110313
$arr_B2[5,0] = 1.96888803690797
$arr_B2[5,1] = 0.4450467218092626
$arr_B2[5,2] = 0.03148727548558483
$arr_B2[5,3] = 0.07781877738676979
$arr_B2[5,4] = 3.305401195893268
$arr_B2[6,0] = 2.116508977307944
$arr_B2[6,1] = 0.4863506388054191
$arr_B2[6,2] = 0.03178439813079592
$arr_B2[6,3] = 0.07755938231196069
$arr_B2[6,4] = 3.356522875535859
$arr_B2[7,0] = 2.417154895887904
$arr_B2[7,1] = 0.5699344432966882
$arr_B2[7,2] = 0.03235310807239777
$arr_B2[7,3] = 0.0772278867995162
$arr_B2[7,4] = 3.47196630502637
$arr_B2[8,0] = 2.645058048405929
$arr_B2[8,1] = 0.6329864294448271
$arr_B2[8,2] = 0.03276369182073147
$arr_B2[8,3] = 0.07709181703292067
$arr_B2[8,4] = 3.566223261354537
$arr_B2[9,0] = 2.750309852434327
$arr_B2[9,1] = 0.6620420758729892
$arr_B2[9,2] = 0.03294918205931019
$arr_B2[9,3] = 0.07705314888689152
$arr_B2[9,4] = 3.611193115692458
$arr_B2[10,0] = 2.790396332125738
$arr_B2[10,1] = 0.6730994220952766
$arr_B2[10,2] = 0.03301926012349909
$arr_B2[10,3] = 0.07704183779620877
$arr_B2[10,4] = 3.628525898786961
$arr_B2[11,0] = 2.781752716863991
$arr_B2[11,1] = 0.6707155805862612
$arr_B2[11,2] = 0.03300417445797876
$arr_B2[11,3] = 0.07704412580201137
$arr_B2[11,4] = 3.624779421606263
$arr_B2[12,0] = 2.753603165805316
$arr_B2[12,1] = 0.6629506704398977
$arr_B2[12,2] = 0.03295495056517339
$arr_B2[12,3] = 0.07705215159422707
$arr_B2[12,4] = 3.612612991494842
$arr_B2[13,0] = 2.736390785550498
$arr_B2[13,1] = 0.658201578879698
$arr_B2[13,2] = 0.03292477893222312
$arr_B2[13,3] = 0.07705750125075816
$arr_B2[13,4] = 3.605200335126199
$arr_B2[14,0] = 2.63821159918723
$arr_B2[14,1] = 0.6310951623417509
$arr_B2[14,2] = 0.03275154557570481
$arr_B2[14,3] = 0.07709481078586045
$arr_B2[14,4] = 3.563326713418974
$arr_B2[15,0] = 2.578388043800487
$arr_B2[15,1] = 0.6145624664036404
$arr_B2[15,2] = 0.03264495850156912
$arr_B2[15,3] = 0.07712364253583104
$arr_B2[15,4] = 3.538176407387112
$arr_B2[16,0] = 2.544127413066064
$arr_B2[16,1] = 0.6050883445625459
$arr_B2[16,2] = 0.03258353001894854
$arr_B2[16,3] = 0.07714241289026802
$arr_B2[16,4] = 3.523907340955134
$arr_B2[17,0] = 2.532552730047826
$arr_B2[17,1] = 0.6018865599617698
$arr_B2[17,2] = 0.0325627097052319
$arr_B2[17,3] = 0.07714914411894114
$arr_B2[17,4] = 3.519109781418706
$arr_B2[18,0] = 2.58474099337559
$arr_B2[18,1] = 0.6163187658636957
$arr_B2[18,2] = 0.03265631739050612
$arr_B2[18,3] = 0.07712034705570048
$arr_B2[18,4] = 3.540833319601859
$arr_B2[19,0] = 2.761865106983635
$arr_B2[19,1] = 0.6652299230233325
$arr_B2[19,2] = 0.03296941306427748
$arr_B2[19,3] = 0.07704970387245602
$arr_B2[19,4] = 3.616178302189354
$arr_B2[20,0] = 2.878967492596189
$arr_B2[20,1] = 0.6975149274951491
$arr_B2[20,2] = 0.03317309921838252
$arr_B2[20,3] = 0.07702295022560257
$arr_B2[20,4] = 3.667192022681007
$arr_B2[21,0] = 2.816343964809278
$arr_B2[21,1] = 0.6802543229737239
$arr_B2[21,2] = 0.03306446673221153
$arr_B2[21,3] = 0.07703545546358725
$arr_B2[21,4] = 3.639801994878383
$arr_B2[22,0] = 2.581868411683615
$arr_B2[22,1] = 0.615524647186021
$arr_B2[22,2] = 0.03265118250626386
$arr_B2[22,3] = 0.07712183010590046
$arr_B2[22,4] = 3.539631537438652
$arr_B2[23,0] = 2.334605377797743
$arr_B2[23,1] = 0.5470391175137479
$arr_B2[23,2] = 0.03220062189915218
$arr_B2[23,3] = 0.07729865296212957
$arr_B2[23,4] = 3.439090955950121
$ws.Range("B2:F25").Value = $arr_B2

$arr_J2 = New-Object "object[,]" 24,1
$arr_J2[0,0] = 0.1480077482287534
$arr_J2[1,0] = 0.1480519533643516
$arr_J2[2,0] = 0.1481383786381478
$arr_J2[3,0] = 0.14818848510156
$arr_J2[4,0] = 0.1481977038444953
$arr_J2[5,0] = 0.1481389941378133
$arr_J2[6,0] = 0.1480106735608331
$arr_J2[7,0] = 0.1482304800102483
$arr_J2[8,0] = 0.1486810823758731
$arr_J2[9,0] = 0.148949248341431
$arr_J2[10,0] = 0.1490599119669369
$arr_J2[11,0] = 0.1490356727096014
$arr_J2[12,0] = 0.1489581698559732
$arr_J2[13,0] = 0.1489118850425797
$arr_J2[14,0] = 0.148664831003579
$arr_J2[15,0] = 0.148529474878579
$arr_J2[16,0] = 0.1484575670596797
$arr_J2[17,0] = 0.1484342405717172
$arr_J2[18,0] = 0.1485432682307177
$arr_J2[19,0] = 0.1489806867076666
$arr_J2[20,0] = 0.1493197083364493
$arr_J2[21,0] = 0.1491338936184832
$arr_J2[22,0] = 0.1485370138492428
$arr_J2[23,0] = 0.1481203722278224
$ws.Range("J2:J25").Value = $arr_J2

$arr_M2 = New-Object "object[,]" 24,2
$arr_M2[0,0] = 0.6361323976315276
$arr_M2[0,1] = 2.230108892598381
$arr_M2[1,0] = 0.6121101459420188
$arr_M2[1,1] = 2.242756627550861
$arr_M2[2,0] = 0.5977641422963131
$arr_M2[2,1] = 2.251219209281018
$arr_M2[3,0] = 0.5920190565311003
$arr_M2[3,1] = 2.254842311893064
$arr_M2[4,0] = 0.5910711768302335
$arr_M2[4,1] = 2.255454446561266
$arr_M2[5,0] = 0.5976862536554748
$arr_M2[5,1] = 2.251267365937139
$arr_M2[6,0] = 0.6277654121183645
$arr_M2[6,1] = 2.234324789872332
$arr_M2[7,0] = 0.6899810506752573
$arr_M2[7,1] = 2.206660632456362
$arr_M2[8,0] = 0.7377037233240458
$arr_M2[8,1] = 2.189768822909784
$arr_M2[9,0] = 0.7598610339173177
$arr_M2[9,1] = 2.18283930610842
$arr_M2[10,0] = 0.7683165694784435
$arr_M2[10,1] = 2.180324608233406
$arr_M2[11,0] = 0.7664926158031307
$arr_M2[11,1] = 2.18086131606114
$arr_M2[12,0] = 0.7605553687665605
$arr_M2[12,1] = 2.182630223804182
$arr_M2[13,0] = 0.7569271226254273
$arr_M2[13,1] = 2.183727997933005
$arr_M2[14,0] = 0.7362647588057314
$arr_M2[14,1] = 2.190236940317135
$arr_M2[15,0] = 0.72370425237159
$arr_M2[15,1] = 2.194423844988492
$arr_M2[16,0] = 0.7165219138568091
$arr_M2[16,1] = 2.196903022169579
$arr_M2[17,0] = 0.7140973152154686
$arr_M2[17,1] = 2.197754594449407
$arr_M2[18,0] = 0.7250369747450947
$arr_M2[18,1] = 2.193970789588292
$arr_M2[19,0] = 0.7622975118964774
$arr_M2[19,1] = 2.182107677955386
$arr_M2[20,0] = 0.7870288449404228
$arr_M2[20,1] = 2.174992280748953
$arr_M2[21,0] = 0.7737943420799382
$arr_M2[21,1] = 2.178731259011727
$arr_M2[22,0] = 0.724434329953894
$arr_M2[22,1] = 2.194175391516694
$arr_M2[23,0] = 0.672799789497752
$arr_M2[23,1] = 2.213544820171904
$ws.Range("M2:N25").Value = $arr_M2
